$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Harvard case classification" column pair (BP/BQ) is relabeled: the column that used
# to hold the "average_doctor" figures now holds "average_doctor_old", and the previous
# "average_doctor_old" column now holds the "average_doctor" figures (values recomputed).
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# New "average_doctor_old" values (previously in BQ, now in BP)
$ws.Range("BP4").Value = 0.65
$ws.Range("BP5").Value = 0.456
$ws.Range("BP6").Value = 0.532
$ws.Range("BP7").Value = 0.483
$ws.Range("BP8").Value = 0.576
$ws.Range("BP9").Value = 0.5600000000000001
$ws.Range("BP10").Value = 0.667
$ws.Range("BP11").Value = 0.667
$ws.Range("BP12").Value = 1.204
$ws.Range("BP13").Value = 0.84

# New "average_doctor" values (previously in BP, now in BQ)
$ws.Range("BQ4").Value = 0.656
$ws.Range("BQ5").Value = 0.454
$ws.Range("BQ6").Value = 0.533
$ws.Range("BQ7").Value = 0.482
$ws.Range("BQ8").Value = 0.586
$ws.Range("BQ9").Value = 0.5590000000000001
$ws.Range("BQ10").Value = 0.6929999999999999
$ws.Range("BQ11").Value = 0.6929999999999999
$ws.Range("BQ12").Value = 1.266
$ws.Range("BQ13").Value = 0.786

# Updated average/variance/std-dev figures for the recalculated "_old" columns
$ws.Range("E4").Value = 0.365
$ws.Range("E5").Value = 0.487
$ws.Range("E6").Value = 0.417
$ws.Range("E7").Value = 0.456
$ws.Range("E8").Value = 0.513
$ws.Range("E9").Value = 0.44
$ws.Range("E10").Value = 0.5600000000000001
$ws.Range("E11").Value = 0.58
$ws.Range("E12").Value = 1.414
$ws.Range("E13").Value = 1.738
$ws.Range("F4").Value = 0.077
$ws.Range("F5").Value = 0.104
$ws.Range("F8").Value = 0.131
$ws.Range("F9").Value = 0.246
$ws.Range("F10").Value = 0.246
$ws.Range("F11").Value = 0.244
$ws.Range("F12").Value = 0.656
$ws.Range("F13").Value = 0.9
$ws.Range("G4").Value = 0.278
$ws.Range("G5").Value = 0.323
$ws.Range("G8").Value = 0.362
$ws.Range("G9").Value = 0.496
$ws.Range("G10").Value = 0.496
$ws.Range("G11").Value = 0.494
$ws.Range("G12").Value = 0.8100000000000001
$ws.Range("G13").Value = 0.949
$ws.Range("N4").Value = 0.384
$ws.Range("N5").Value = 0.752
$ws.Range("N6").Value = 0.508
$ws.Range("N7").Value = 0.631
$ws.Range("N8").Value = 0.747
$ws.Range("N9").Value = 0.62
$ws.Range("N10").Value = 0.82
$ws.Range("N11").Value = 0.86
$ws.Range("N12").Value = 1.667
$ws.Range("N13").Value = 2.378
$ws.Range("O4").Value = 0.058
$ws.Range("O5").Value = 0.08599999999999999
$ws.Range("O8").Value = 0.07199999999999999
$ws.Range("O9").Value = 0.236
$ws.Range("O10").Value = 0.148
$ws.Range("O11").Value = 0.12
$ws.Range("O12").Value = 1.644
$ws.Range("O13").Value = 1.098
$ws.Range("P4").Value = 0.241
$ws.Range("P5").Value = 0.293
$ws.Range("P8").Value = 0.268
$ws.Range("P9").Value = 0.485
$ws.Range("P10").Value = 0.384
$ws.Range("P11").Value = 0.347
$ws.Range("P12").Value = 1.282
$ws.Range("P13").Value = 1.048
$ws.Range("W4").Value = 0.217
$ws.Range("W5").Value = 0.225
$ws.Range("W6").Value = 0.221
$ws.Range("W7").Value = 0.223
$ws.Range("W8").Value = 0.217
$ws.Range("W9").Value = 0.12
$ws.Range("W10").Value = 0.26
$ws.Range("W11").Value = 0.26
$ws.Range("W12").Value = 1.846
$ws.Range("W13").Value = 1.09
$ws.Range("X4").Value = 0.101
$ws.Range("X5").Value = 0.113
$ws.Range("X8").Value = 0.105
$ws.Range("X9").Value = 0.106
$ws.Range("X10").Value = 0.192
$ws.Range("X11").Value = 0.192
$ws.Range("X12").Value = 0.746
$ws.Range("X13").Value = 0.186
$ws.Range("Y4").Value = 0.318
$ws.Range("Y5").Value = 0.337
$ws.Range("Y8").Value = 0.325
$ws.Range("Y9").Value = 0.325
$ws.Range("Y10").Value = 0.439
$ws.Range("Y11").Value = 0.439
$ws.Range("Y12").Value = 0.863
$ws.Range("Y13").Value = 0.431
$ws.Range("AI4").Value = 0.197
$ws.Range("AI5").Value = 0.235
$ws.Range("AI6").Value = 0.214
$ws.Range("AI7").Value = 0.226
$ws.Range("AI8").Value = 0.215
$ws.Range("AI9").Value = 0.12
$ws.Range("AI10").Value = 0.24
$ws.Range("AI11").Value = 0.24
$ws.Range("AI12").Value = 1.917
$ws.Range("AI13").Value = 1.383
$ws.Range("AJ4").Value = 0.064
$ws.Range("AJ5").Value = 0.093
$ws.Range("AJ8").Value = 0.092
$ws.Range("AJ9").Value = 0.106
$ws.Range("AJ10").Value = 0.182
$ws.Range("AJ11").Value = 0.182
$ws.Range("AJ12").Value = 0.91
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK4").Value = 0.254
$ws.Range("AK5").Value = 0.305
$ws.Range("AK8").Value = 0.303
$ws.Range("AK9").Value = 0.325
$ws.Range("AK10").Value = 0.427
$ws.Range("AK11").Value = 0.427
$ws.Range("AK12").Value = 0.954
$ws.Range("AK13").Value = 0.633
$ws.Range("AU4").Value = 0.144
$ws.Range("AU5").Value = 0.294
$ws.Range("AU6").Value = 0.193
$ws.Range("AU7").Value = 0.243
$ws.Range("AU8").Value = 0.233
$ws.Range("AU10").Value = 0.22
$ws.Range("AU11").Value = 0.32
$ws.Range("AU12").Value = 2.562
$ws.Range("AU13").Value = 2.481
$ws.Range("AV4").Value = 0.026
$ws.Range("AV5").Value = 0.099
$ws.Range("AV8").Value = 0.076
$ws.Range("AV10").Value = 0.172
$ws.Range("AV11").Value = 0.218
$ws.Range("AV12").Value = 1.746
$ws.Range("AV13").Value = 1.334
$ws.Range("AW4").Value = 0.162
$ws.Range("AW5").Value = 0.314
$ws.Range("AW8").Value = 0.275
$ws.Range("AW10").Value = 0.414
$ws.Range("AW11").Value = 0.466
$ws.Range("AW12").Value = 1.321
$ws.Range("AW13").Value = 1.155
$ws.Range("BA4").Value = 1.949
$ws.Range("BA5").Value = 1.367
$ws.Range("BA6").Value = 1.595
$ws.Range("BA7").Value = 1.448
$ws.Range("BA8").Value = 1.727
$ws.Range("BA9").Value = 1.68
$ws.Range("BA10").Value = 2
$ws.Range("BA11").Value = 2
$ws.Range("BA12").Value = 3.611
$ws.Range("BA13").Value = 2.519
$ws.Range("BB4").Value = 0.169
$ws.Range("BB5").Value = 0.08599999999999999
$ws.Range("BB8").Value = 0.137
$ws.Range("BB9").Value = 0.246
$ws.Range("BB12").Value = 0.32
$ws.Range("BB13").Value = 0.313
$ws.Range("BC4").Value = 0.411
$ws.Range("BC5").Value = 0.293
$ws.Range("BC8").Value = 0.37
$ws.Range("BC9").Value = 0.496
$ws.Range("BC12").Value = 0.5659999999999999
$ws.Range("BC13").Value = 0.5590000000000001
$ws.Range("BG4").Value = 0.73
$ws.Range("BG5").Value = 0.41
$ws.Range("BG6").Value = 0.525
$ws.Range("BG7").Value = 0.449
$ws.Range("BG8").Value = 0.5649999999999999
$ws.Range("BG9").Value = 0.6
$ws.Range("BG10").Value = 0.64
$ws.Range("BG11").Value = 0.64
$ws.Range("BG12").Value = 1.062
$ws.Range("BG13").Value = 0.621
$ws.Range("BH4").Value = 0.142
$ws.Range("BH8").Value = 0.105
$ws.Range("BH9").Value = 0.24
$ws.Range("BH10").Value = 0.23
$ws.Range("BH11").Value = 0.23
$ws.Range("BH12").Value = 0.059
$ws.Range("BH13").Value = 0.08699999999999999
$ws.Range("BI4").Value = 0.377
$ws.Range("BI5").Value = 0.229
$ws.Range("BI8").Value = 0.324
$ws.Range("BI9").Value = 0.49
$ws.Range("BI10").Value = 0.48
$ws.Range("BI11").Value = 0.48
$ws.Range("BI12").Value = 0.242
$ws.Range("BI13").Value = 0.294
$ws.Range("BM4").Value = 0.6820000000000001
$ws.Range("BM5").Value = 0.574
$ws.Range("BM6").Value = 0.623
$ws.Range("BM7").Value = 0.593
$ws.Range("BM8").Value = 0.702
$ws.Range("BM9").Value = 0.64
$ws.Range("BM10").Value = 0.86
$ws.Range("BM11").Value = 0.86
$ws.Range("BM12").Value = 1.349
$ws.Range("BM13").Value = 0.999
$ws.Range("BN4").Value = 0.091
$ws.Range("BN8").Value = 0.07199999999999999
$ws.Range("BN9").Value = 0.23
$ws.Range("BN10").Value = 0.12
$ws.Range("BN11").Value = 0.12
$ws.Range("BN12").Value = 0.413
$ws.Range("BN13").Value = 0.363
$ws.Range("BO4").Value = 0.302
$ws.Range("BO8").Value = 0.268
$ws.Range("BO9").Value = 0.48
$ws.Range("BO10").Value = 0.347
$ws.Range("BO11").Value = 0.347
$ws.Range("BO12").Value = 0.643
$ws.Range("BO13").Value = 0.603
